$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row cells: "_old" -> "_FV2410", "_new" -> "_FV2504"
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value2
    if (($v -ne $null) -and ($v -is [string])) {
        if ($v.EndsWith("_old")) {
            $cell.Value2 = $v.Substring(0, $v.Length - 4) + "_FV2410"
        } elseif ($v.EndsWith("_new")) {
            $cell.Value2 = $v.Substring(0, $v.Length - 4) + "_FV2504"
        }
    }
}

# Freeze header row (pane split after row 1)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the used range into an Excel Table ("Table1") with a header row
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U65"), $null, 1)
$tbl.Name = "Table1"
